# Insert a new data row at row 121 (pushing the existing rows 121-155 down
# to 122-156), then populate the new row with the new weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(121).Insert()

$ws.Range("A121").Value = 4
$ws.Range("B121").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C121").Value = 'Los Lagos'
$ws.Range("D121").Value = 44508
$ws.Range("E121").Value = 10
$ws.Range("F121").Value = 'Fruta'
$ws.Range("G121").Value = 100108
$ws.Range("H121").Value = 'Tropicales y subtropicales'
$ws.Range("I121").Value = 100108005
$ws.Range("J121").Value = 'Piña'
$ws.Range("K121").Value = 'Caramelo'
$ws.Range("L121").Value = 'Segunda'
$ws.Range("M121").Value = 80
$ws.Range("N121").Value = 21000
$ws.Range("O121").Value = 22000
$ws.Range("P121").Value = 21500
$ws.Range("Q121").Value = '$/caja 14 unidades'
$ws.Range("R121").Value = 'Ecuador'
$ws.Range("S121").Value = 1536
$ws.Range("T121").Value = 14
